$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source workbook ships with an empty styleSheet (no cellXfs entries),
# which makes any direct value write crash the engine's style lookup.
# Touching .Style first forces a default style entry to be created so
# subsequent writes are safe.
$ws.Range("A1").Style = "Normal"

# Shift the existing two data rows down to make room for a header row.
$ws.Rows("1:1").Insert()

# Write the new header row.
$ws.Range("A1").Value = "a"
$ws.Range("B1").Value = "b"
$ws.Range("C1").Value = "c"
